$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 62-91 (Simulation Controls section) down by 30 rows to make room
# for the new APSIM soil parameter blocks.
$ws.Rows("62:91").Insert(-4121)

# ---------------------------------------------------------------------------
# Block 1: Soil water content, air dry (row 62-63)
# ---------------------------------------------------------------------------
$ws.Range("A62").Value = "!"
$ws.Range("C62").Value = "Soil water content, air dry"
$ws.Range("D62").Value = "Multiply function"
$ws.Range("E62").Value = "Soil water, lower limit"
$ws.Range("F62").Value = "multiplier"

$ws.Range("A63").Value = "&"
$ws.Range("B63").Value = "FILL"
$ws.Range("C63").Value = "SLADR"
$ws.Range("D63").Value = "MULTIPLY()"
$ws.Range("E63").Value = "`$SLLL"
$ws.Range("F63").Value = 0.4

# ---------------------------------------------------------------------------
# Block 2: APSIM Evaporative diffusion constant (row 65-66)
# ---------------------------------------------------------------------------
$ws.Range("A65").Value = "!"
$ws.Range("C65").Value = "APSIM Evaporative diffusion constant"
$ws.Range("D65").Value = "value"

$ws.Range("A66").Value = "&"
$ws.Range("B66").Value = "FILL"
$ws.Range("C66").Value = "DiffusConst"
$ws.Range("D66").Value = 40

# ---------------------------------------------------------------------------
# Block 3: APSIM Evaporative diffusion slope (row 68-69)
# ---------------------------------------------------------------------------
$ws.Range("A68").Value = "!"
$ws.Range("C68").Value = "APSIM Evaporative diffusion slope"
$ws.Range("D68").Value = "value"

$ws.Range("A69").Value = "&"
$ws.Range("B69").Value = "FILL"
$ws.Range("C69").Value = "DiffusSlope"
$ws.Range("D69").Value = 16

# ---------------------------------------------------------------------------
# Block 4: Apsim Cona soil evaporation (row 71-72)
# ---------------------------------------------------------------------------
$ws.Range("A71").Value = "!"
$ws.Range("C71").Value = "Apsim Cona soil evaporation"
$ws.Range("D71").Value = "value"

$ws.Range("A72").Value = "&"
$ws.Range("B72").Value = "FILL"
$ws.Range("C72").Value = "CONA"
$ws.Range("D72").Value = 3.5

# ---------------------------------------------------------------------------
# Block 5: Summer soil evaporation, upper limit for stage 1 (row 74-75)
# ---------------------------------------------------------------------------
$ws.Range("A74").Value = "!"
$ws.Range("C74").Value = "Summer soil evaporation, upper limit for stage 1"
$ws.Range("D74").Value = "value"

$ws.Range("A75").Value = "&"
$ws.Range("B75").Value = "FILL"
$ws.Range("C75").Value = "SummerU"
$ws.Range("D75").Value = 6

# ---------------------------------------------------------------------------
# Block 6: Winter soil evaporation, upper limit for stage 1 (row 77-78)
# ---------------------------------------------------------------------------
$ws.Range("A77").Value = "!"
$ws.Range("C77").Value = "Winter soil evaporation, upper limit for stage 1"
$ws.Range("D77").Value = "value"

$ws.Range("A78").Value = "&"
$ws.Range("B78").Value = "FILL"
$ws.Range("C78").Value = "WinterU"
$ws.Range("D78").Value = 6

# ---------------------------------------------------------------------------
# Block 7: Date to switch to SummerU (row 80-81)
# ---------------------------------------------------------------------------
$ws.Range("A80").Value = "!"
$ws.Range("C80").Value = "Date to switch to SummerU "
$ws.Range("D80").Value = "value"

$ws.Range("A81").Value = "&"
$ws.Range("B81").Value = "FILL"
$ws.Range("C81").Value = "SummerDate__soil"
$ws.Range("D81").Value = "03-21"

# ---------------------------------------------------------------------------
# Block 8: Date to switch to WinterU (row 83-84)
# ---------------------------------------------------------------------------
$ws.Range("A83").Value = "!"
$ws.Range("C83").Value = "Date to switch to WinterU"
$ws.Range("D83").Value = "value"

$ws.Range("A84").Value = "&"
$ws.Range("B84").Value = "FILL"
$ws.Range("C84").Value = "WinterDate__soil"
$ws.Range("D84").Value = "09-21"

# ---------------------------------------------------------------------------
# Block 9: APSIM plant specific water uptake coefficient (row 86-87)
# ---------------------------------------------------------------------------
$ws.Range("A86").Value = "!"
$ws.Range("C86").Value = "APSIM plant specific water uptake coefficient"
$ws.Range("D86").Value = "value"

$ws.Range("A87").Value = "&"
$ws.Range("B87").Value = "FILL"
$ws.Range("C87").Value = "apsim_kl"
$ws.Range("D87").Value = 0.06

# ---------------------------------------------------------------------------
# Block 10: Biologically active soil organic carbon by layer (row 89-90)
# ---------------------------------------------------------------------------
$ws.Range("A89").Value = "!"
$ws.Range("C89").Value = "Biologically active soil organic carbon by layer"
$ws.Range("D89").Value = "value"

$ws.Range("A90").Value = "&"
$ws.Range("B90").Value = "FILL"
$ws.Range("C90").Value = "slacc"
$ws.Range("D90").Value = 0.02

$wb.Save()
